# Update "想去人数" (interested count) values in F column on the
# "展览" and "全部类型" worksheets, mirroring the latest site scrape.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 6645
    6  = 2032
    7  = 1545
    8  = 308
    9  = 1013
    10 = 430
    12 = 5637
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
